# Update existing rows 2-5 with new timestamps / run_time values,
# then delete rows 6-11 (only 4 data rows remain after the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "2025-07-01 23:52:34"
$ws.Range("H2").Value = 0.2879897669990896
$ws.Range("I2").Value = 0.004799829449984827

# Row 3
$ws.Range("A3").Value = "2025-07-01 23:52:35"
$ws.Range("H3").Value = 0.2778664880024735
$ws.Range("I3").Value = 0.004631108133374558

# Row 4
$ws.Range("A4").Value = "2025-07-01 23:52:39"
$ws.Range("H4").Value = 4.519339108999702
$ws.Range("I4").Value = 0.07532231848332836

# Row 5
$ws.Range("A5").Value = "2025-07-01 23:52:45"
$ws.Range("H5").Value = 5.357828082000196
$ws.Range("I5").Value = 0.08929713470000328

# Remove rows 6 through 11 (trailing data no longer present).
$ws.Range("A6:I11").Delete()
